$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2470.6875
$ws.Range("I15").Value = 2470.6875
$ws.Range("K15").Value = 7412.0625
$ws.Range("M15").Value = -7243.0625
$ws.Range("H17").Value = 10001984
$ws.Range("J17").Value = 10001984
$ws.Range("L17").Value = 30005952
$ws.Range("N17").Value = -30006288
$ws.Range("H51").Value = 15402.25
$ws.Range("J51").Value = 9805
$ws.Range("L51").Value = 9805
$ws.Range("N51").Value = -10773
$ws.Range("H127").Value = 11548.25
$ws.Range("I127").Value = 11548.25
$ws.Range("K127").Value = 34644.75
$ws.Range("M127").Value = -29684.75
$ws.Range("H137").Value = 1733.2142
$ws.Range("I137").Value = 1712.6923
$ws.Range("K137").Value = 5138.0769
$ws.Range("M137").Value = -2588.0769

# --- ARM sheet updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2998.8462
$ws.Range("I32").Value = 2492.7192
$ws.Range("K32").Value = 2492.7192
$ws.Range("M32").Value = -2205.7192
$ws.Range("H61").Value = 3973109.2
$ws.Range("I61").Value = 4633213.5
$ws.Range("K61").Value = 4633213.5
$ws.Range("M61").Value = -4633001.5
$ws.Range("H74").Value = 3350.2856
$ws.Range("I74").Value = 2685.75
$ws.Range("K74").Value = 2685.75
$ws.Range("M74").Value = -1811.75
$ws.Range("H77").Value = 3350.2856
$ws.Range("I77").Value = 2685.75
$ws.Range("K77").Value = 13428.75
$ws.Range("M77").Value = -9060.75
$ws.Range("H88").Value = 1281.1177
$ws.Range("I88").Value = 1240.8572
$ws.Range("J88").Value = 1309.3
$ws.Range("K88").Value = 1240.8572
$ws.Range("L88").Value = 1309.3
$ws.Range("M88").Value = -834.8571999999999
$ws.Range("N88").Value = -2121.3
$ws.Range("H91").Value = 1281.1177
$ws.Range("I91").Value = 1240.8572
$ws.Range("J91").Value = 1309.3
$ws.Range("K91").Value = 1240.8572
$ws.Range("L91").Value = 1309.3
$ws.Range("M91").Value = 163.1428000000001
$ws.Range("N91").Value = -4117.3
$ws.Range("H95").Value = 21666
$ws.Range("J95").Value = 21666
$ws.Range("L95").Value = 21666
$ws.Range("N95").Value = -27158
$ws.Range("H104").Value = 14999.5
$ws.Range("J104").Value = 14999.5
$ws.Range("L104").Value = 14999.5
$ws.Range("N104").Value = -21987.5
$ws.Range("H136").Value = 3973109.2
$ws.Range("I136").Value = 4633213.5
$ws.Range("K136").Value = 13899640.5
$ws.Range("M136").Value = -13897090.5
$ws.Range("H141").Value = 98429
$ws.Range("J141").Value = 98429
$ws.Range("L141").Value = 98429
$ws.Range("N141").Value = -108789

# --- BSM sheet updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5022.294
$ws.Range("I99").Value = 3788.4
$ws.Range("K99").Value = 3788.4
$ws.Range("M99").Value = -2290.4
$ws.Range("H107").Value = 3891.3635
$ws.Range("I107").Value = 4034.5293
$ws.Range("J107").Value = 3404.6
$ws.Range("K107").Value = 4034.5293
$ws.Range("L107").Value = 3404.6
$ws.Range("M107").Value = -2114.5293
$ws.Range("N107").Value = -7244.6
$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680
$ws.Range("H116").Value = 74963.336
$ws.Range("I116").Value = 145000
$ws.Range("J116").Value = 39945
$ws.Range("K116").Value = 145000
$ws.Range("L116").Value = 39945
$ws.Range("M116").Value = -140411
$ws.Range("N116").Value = -49123
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 5183.7915
$ws.Range("I134").Value = 5072.2173
$ws.Range("K134").Value = 15216.6519
$ws.Range("M134").Value = -12681.6519

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2082.7827
$ws.Range("I16").Value = 1807.7333
$ws.Range("J16").Value = 2598.5
$ws.Range("K16").Value = 1807.7333
$ws.Range("L16").Value = 2598.5
$ws.Range("M16").Value = -1520.7333
$ws.Range("N16").Value = -3172.5
$ws.Range("H22").Value = 210.66667
$ws.Range("I22").Value = 175.1
$ws.Range("J22").Value = 281.8
$ws.Range("K22").Value = 175.1
$ws.Range("L22").Value = 281.8
$ws.Range("M22").Value = 174.9
$ws.Range("N22").Value = -981.8
$ws.Range("H105").Value = 1137.4166
$ws.Range("I105").Value = 1039
$ws.Range("J105").Value = 1432.6666
$ws.Range("K105").Value = 1039
$ws.Range("L105").Value = 1432.6666
$ws.Range("M105").Value = 708
$ws.Range("N105").Value = -4926.6666
$ws.Range("H107").Value = 3762.3845
$ws.Range("I107").Value = 8237
$ws.Range("J107").Value = 2420
$ws.Range("K107").Value = 8237
$ws.Range("L107").Value = 2420
$ws.Range("M107").Value = -6317
$ws.Range("N107").Value = -6260
$ws.Range("H113").Value = 2082.7827
$ws.Range("I113").Value = 1807.7333
$ws.Range("J113").Value = 2598.5
$ws.Range("K113").Value = 1807.7333
$ws.Range("L113").Value = 2598.5
$ws.Range("M113").Value = 362.2666999999999
$ws.Range("N113").Value = -6938.5
$ws.Range("H134").Value = 7445.107
$ws.Range("I134").Value = 4263.75
$ws.Range("K134").Value = 12791.25
$ws.Range("M134").Value = -10256.25

# --- CUL sheet updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6082.5
$ws.Range("I3").Value = 6082.5
$ws.Range("K3").Value = 18247.5
$ws.Range("M3").Value = -18135.5
$ws.Range("H5").Value = 1151.1628
$ws.Range("I5").Value = 1016.7619
$ws.Range("K5").Value = 3050.2857
$ws.Range("M5").Value = -2938.2857
$ws.Range("H18").Value = 1273.8182
$ws.Range("I18").Value = 626.75
$ws.Range("K18").Value = 1880.25
$ws.Range("M18").Value = -1711.25
$ws.Range("H134").Value = 1226.3334
$ws.Range("I134").Value = 1226.3334
$ws.Range("K134").Value = 3679.0002
$ws.Range("M134").Value = 1390.9998
$ws.Range("H135").Value = 1151.1628
$ws.Range("I135").Value = 1016.7619
$ws.Range("K135").Value = 9150.857099999999
$ws.Range("M135").Value = -6615.857099999999

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 39111
$ws.Range("J93").Value = 39111
$ws.Range("L93").Value = 39111
$ws.Range("N93").Value = -42855

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7445.1577
$ws.Range("I7").Value = 8300.846
$ws.Range("K7").Value = 8300.846
$ws.Range("M7").Value = -8188.846
$ws.Range("H122").Value = 2934.6924
$ws.Range("I122").Value = 2375.25
$ws.Range("J122").Value = 3829.8
$ws.Range("K122").Value = 7125.75
$ws.Range("L122").Value = 11489.4
$ws.Range("M122").Value = -4675.75
$ws.Range("N122").Value = -16389.4
$ws.Range("H126").Value = 7445.1577
$ws.Range("I126").Value = 8300.846
$ws.Range("K126").Value = 24902.538
$ws.Range("M126").Value = -22432.538
$ws.Range("H136").Value = 3203.15
$ws.Range("I136").Value = 3004.0557
$ws.Range("K136").Value = 9012.167099999999
$ws.Range("M136").Value = -6462.167099999999

# --- WVR sheet updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 85164.836
$ws.Range("J46").Value = 85164.836
$ws.Range("L46").Value = 85164.836
$ws.Range("N46").Value = -85626.836
$ws.Range("H81").Value = 5392.2856
$ws.Range("J81").Value = 6359.2
$ws.Range("L81").Value = 12718.4
$ws.Range("N81").Value = -14840.4
$ws.Range("H84").Value = 5392.2856
$ws.Range("J84").Value = 6359.2
$ws.Range("L84").Value = 63592
$ws.Range("N84").Value = -74200
$ws.Range("H122").Value = 3905.3447
$ws.Range("I122").Value = 3335.9333
$ws.Range("K122").Value = 10007.7999
$ws.Range("M122").Value = -7557.7999
$ws.Range("H126").Value = 5592
$ws.Range("I126").Value = 2913.4666
$ws.Range("K126").Value = 8740.399800000001
$ws.Range("M126").Value = -6270.399800000001
$ws.Range("H134").Value = 85164.836
$ws.Range("J134").Value = 85164.836
$ws.Range("L134").Value = 255494.508
$ws.Range("N134").Value = -260564.508
$ws.Range("H136").Value = 3359.5625
$ws.Range("I136").Value = 2549.2163
$ws.Range("K136").Value = 7647.6489
$ws.Range("M136").Value = -5097.6489

